$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '93.646.26'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').Value = '3.433.18'
$ws.Range('E3').Value = '  +2.44%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '234.55'
$ws.Range('E5').Value = '  +0.98%  '
$ws.Range('D6').Value = '622.56'
$ws.Range('E6').Value = '  -1.03%  '
$ws.Range('E7').Value = '  +7.12%  '
$ws.Range('D8').Value = '0.395'
$ws.Range('E8').Value = '  +2.18%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').Value = '0.987'
$ws.Range('E10').Value = '  +5.97%  '
$ws.Range('D11').Value = '3.430.23'
$ws.Range('E11').Value = '  +2.45%  '
$ws.Range('D12').Value = '43.88'
$ws.Range('E12').Value = '  +8.89%  '
$ws.Range('E13').Value = '  +2.71%  '
$ws.Range('D14').Value = '6.29'
$ws.Range('E14').Value = '  +5.79%  '
$ws.Range('D15').Value = '93.471.39'
$ws.Range('E15').Value = '  +0.45%  '
$ws.Range('D16').Value = '4.069.94'
$ws.Range('E16').Value = '  +2.49%  '
$ws.Range('D17').Value = '0.0000248'
$ws.Range('E17').Value = '  +2.25%  '
$ws.Range('D18').Value = '8.31'
$ws.Range('E18').Value = '  +4.83%  '
$ws.Range('D19').Value = '3.419.52'
$ws.Range('E19').Value = '  +2.14%  '
$ws.Range('D20').Value = '18.18'
$ws.Range('E20').Value = '  +8.21%  '
$ws.Range('D21').Value = '11.65'
$ws.Range('E21').Value = '  +6.62%  '
$ws.Range('D22').Value = '0.524'
$ws.Range('E22').Value = '  +16.89%  '
$ws.Range('E23').Value = '  +8.32%  '
$ws.Range('D24').Value = '501.05'
$ws.Range('E24').Value = '  +2.18%  '
$ws.Range('D25').Value = '6.74'
$ws.Range('E25').Value = '  +8.01%  '
$ws.Range('D26').Value = '0.0000184'
$ws.Range('E26').Value = '  -0.28%  '
$ws.Range('D27').Value = '87.98'
$ws.Range('E27').Value = '  -1.38%  '
$ws.Range('D28').Value = '12.16'
$ws.Range('E28').Value = '  +6.42%  '
$ws.Range('D29').Value = '3.613.14'
$ws.Range('E29').Value = '  +2.51%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = '0.142'
$ws.Range('E30').Value = '  +8.33%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '11.42'
$ws.Range('E31').Value = '  +1.25%  '
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('D33').Value = '2.74'
$ws.Range('E33').Value = '  +3.98%  '
$ws.Range('D34').Value = '0.983'
$ws.Range('E34').Value = '  -1.56%  '
$ws.Range('D35').Value = '0.178'
$ws.Range('E35').Value = '  +4.49%  '
$ws.Range('E36').Value = '  +6.10%  '
$ws.Range('D37').Value = '29.27'
$ws.Range('E37').Value = '  +3.05%  '
$ws.Range('D38').Value = '569.05'
$ws.Range('E38').Value = '  +9.13%  '
$ws.Range('D39').Value = '7.55'
$ws.Range('E39').Value = '  +1.51%  '
$ws.Range('D40').Value = '1.42'
$ws.Range('E40').Value = '  +2.33%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('E42').Value = '  +1.58%  '
$ws.Range('D43').Value = '0.905'
$ws.Range('E43').Value = '  +3.79%  '
$ws.Range('D44').Value = '23.74'
$ws.Range('E44').Value = '  -1.23%  '
$ws.Range('E45').Value = '  +1.87%  '
$ws.Range('D46').Value = '0.0418'
$ws.Range('E46').Value = '  +6.35%  '
$ws.Range('D47').Value = '3.64'
$ws.Range('E47').Value = '  +2.70%  '
$ws.Range('D48').Value = '5.51'
$ws.Range('E48').Value = '  +1.08%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').Value = '53.25'
$ws.Range('E49').Value = '  +2.26%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '8.16'
$ws.Range('E50').Value = '  +3.66%  '
$ws.Range('D51').Value = '2.12'
$ws.Range('E51').Value = '  -1.32%  '
